$wb = $excel.ActiveWorkbook

# Rename the "Device" sheet to "Apparatus"
$ws = $wb.Worksheets.Item("Device")
$ws.Name = "Apparatus"

# Copy the bold header style already used on A1/A2 onto the A3:C3 header row
$ws.Range("A1").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Update the sheet's descriptive text (device -> apparatus)
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Update the selected range shown when the sheet is active
[void]$ws.Range("A3:C3").Select()
